$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "id" column (old column A). Everything shifts one column to
# the left (B->A, C->B, ... L->K) as a side effect.
$ws.Columns.Item(1).Delete()

# --- Data fixes on the shifted sheet (columns now A:K) ---

# zip_code for row 2 (Finland) and row 3 (Germany) were stored as
# 6-digit numbers; correct them to the real 5-digit zip codes.
$ws.Range("H2").Value = 32323
$ws.Range("H3").Value = 54645

# Row 4 was incomplete/wrong: country was "Poland" with no address and a
# zip code stored as text; fix it up to be a full France record like the
# other rows.
$ws.Range("F4").Value = "France"
$ws.Range("G4").Value = "Address 1 22-333 4/5"
$ws.Range("H4").Value = 70003

# Phone number for row 4 was stored as text; convert it to a real number.
$ws.Range("J4").NumberFormatLocal = "General"
$ws.Range("J4").Value = 123456789

# Reflect where the user ended up after making the edits.
$ws.Range("J5").Select()
